# lipidcane_spearman_-1.xlsx : refresh Spearman correlation figures and
# extend the parameter/element table with the new 2G (cellulosic) unit
# parameters (crude glycerol, pure glycerine, cellulase streams,
# pretreatment & cofermentation rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the existing correlation values for rows 4-14 (element /
#    parameter rows already on the sheet). Labels in column A/B are
#    unchanged -- only the Spearman rho values move.
# ---------------------------------------------------------------------

$rowData = @{
    4  = @{ C = -0.02121219675248787;  E = 0.01133214832528593;   F = -0.01875564372622575;  H = -0.01932827549313102;  I = -0.01894716392588655;  J = 0.01196679858300604 }
    5  = @{ C = -0.01027079695483188;  E = -0.002859720594388823; F = 0.0004260152810406112;  H = 0.001045199177807967;   I = -0.00162593900903756;  J = 0.007470819134243793 }
    6  = @{ C = 0.005565642654625705;  E = -0.008901098468043937; F = -0.002212107448484297;  H = -0.004335209549408381;  I = -0.002647487913899516; J = -0.001052227878661404 }
    7  = @{ C = 0.003830104665204186;  E = -0.0119754650870186;   F = 0.004614277816571112;   H = 0.00527321877092875;    I = -0.0005958812398352496;J = -0.01427097217562373 }
    8  = @{ C = 0.1091581988623279;    E = -0.0251881402075256;   F = 0.9864043686721746;     H = 0.999999703071988;      I = 0.9576457947858317;    J = -0.05232393255469991 }
    9  = @{ C = 0.9547143702525748;    E = 0.02177359546294381;   F = -0.02881871548874862;   H = -0.02724712131388485;   I = -0.02565199811407992;  J = -0.0114884708325753 }
    10 = @{ C = 0.004275035115001404;  E = -0.0302864602194584;   F = 0.01230549783621991;    H = 0.01120683712027348;    I = 0.007491558251662329;  J = 0.00678602765793675 }
    11 = @{ C = -0.002858188050327522; E = 0.01224121306564852;   F = 0.004479560723182428;   H = 0.005149690765987629;   I = 0.008194044615761783;  J = 0.01332652281091606 }
    12 = @{ C = 0.06794655343786213;   E = 0.0009608923584356942; F = 0.02226205701848227;    H = 0.02384258332170333;    I = 0.02000939523237581;   J = -0.003671776165935196 }
    13 = @{ C = 0.1030724379788975;    E = [double]"-3.0021025200841E-05"; F = -0.01293443629337745; H = -0.01320797592031903; I = 0.267955776926231;  J = 0.01693974752792456 }
    14 = @{ C = -0.199130951997238;    E = 0.004154904454196178;  F = 0.02337183309487332;    H = 0.02070701910028076;    I = 0.02607368293094731;   J = 0.003341824777071115 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("C$r").Value = $vals.C
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
}

# ---------------------------------------------------------------------
# 2) Merge the "Element" column cells that belong together FIRST
#    (matching the existing A4:A8 / A12:A14 merge pattern), before any
#    formatting/values are applied, so the merge operation does not
#    re-split the shared border style into new per-edge styles.
# ---------------------------------------------------------------------

$ws.Range("A17:A18").Merge()
$ws.Range("A20:A21").Merge()
$ws.Range("A22:A23").Merge()

# ---------------------------------------------------------------------
# 3) Append the new rows (15-23) describing the glycerol/cellulase
#    streams and the pretreatment/cofermentation unit parameters.
#    Formats are copied from already-styled rows so no new cell styles
#    get created.
# ---------------------------------------------------------------------

function Set-StyledValue($cellRef, $value, $templateRef) {
    $ws.Range($cellRef).Value = $value
    $ws.Range($templateRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

$newRows = @(
    @{ Row = 15; A = "Stream-crude glycerol";              B = "Price [USD/kg]";                    C = -0.0135902301756092;   E = 0.009334619989384799;   F = -0.001652742786109711; H = 0.002053032658121306;   I = 0.004215266088610644;   J = -0.0145836245990908 }
    @{ Row = 16; A = "Stream-pure glycerine";               B = "Price [USD/kg]";                    C = -0.005417894232715768; E = 0.00741367450454698;    F = 0.02319113267164531;   H = 0.02683620145744805;    I = 0.02145504949820198;    J = -0.003247869444108686 }
    @{ Row = 17; A = "Stream-cellulase";                    B = "Price [USD/kg]";                    C = 0.00874159177366367;   E = -0.006627397321095892;  F = 0.04133105925324237;   H = 0.04119099668763986;    I = 0.04796831395073255;    J = 0.02058471401919945 }
    @{ Row = 18; A = $null;                                 B = "Cellulase loading [wt. % cellulose]";C = 0.02785799487431979;   E = 0.002684852171394086;   F = 0.01331132270845291;   H = 0.01336800773472031;    I = 0.01683920947356838;    J = 0.003519364899805935 }
    @{ Row = 19; A = "Pretreatment reactor system";         B = "Base cost [million USD]";           C = 0.0113997099599884;    E = 0.01252524731700989;    F = -0.001650367362014694; H = 0.002070557554822302;   I = 0.004456985650279426;   J = -0.02312707169414742 }
    @{ Row = 20; A = "Pretreatment and saccharification";   B = "Glucose yield [%]";                 C = 0.009629132353165291;  E = 0.004898395491935819;   F = 0.0001707341828293673; H = -0.0006025162801006511; I = -0.001009648360385934;  J = -0.005255548383826097 }
    @{ Row = 21; A = $null;                                 B = "Xylose yield [%]";                  C = 0.02412859286914371;   E = -0.01265854485834179;   F = 0.02518422225536888;   H = 0.02411183280447331;    I = 0.02411104560444182;    J = 0.01942925905578899 }
    @{ Row = 22; A = "Cofermenation";                       B = "Glucose to ethanol yield [%]";      C = 0.01149286020371441;   E = 0.02069680690787228;    F = [double]"-3.233904129356165E-05"; H = -0.0007787600951504037; I = -0.0009385170615406823; J = -0.006837923689731783 }
    @{ Row = 23; A = $null;                                 B = "Xylose to ethanol yield [%]";       C = -0.008647588377903534; E = -0.008701865820074631;  F = -0.007082371867294874; H = -0.008510149204405966;  I = -0.01020257877610315;   J = -0.002345641358777564 }
)

foreach ($row in $newRows) {
    $r = $row.Row

    if ($row.A -ne $null) {
        Set-StyledValue "A$r" $row.A "A4"
    } else {
        Set-StyledValue "A$r" "" "A5"
        $ws.Range("A$r").Value = $null
    }

    Set-StyledValue "B$r" $row.B "B4"

    $ws.Range("C$r").Value = $row.C
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
}
